$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns C:F (runs, balls, fours, sixes) for the data rows (2-8) are
# rotated down by one row: the last row's stats become the new first
# row's stats, and every other row shifts down by one.
#
# Force the cells to stay text-typed (the sheet stores numeric-looking
# values as text) before writing the new values.
$ws.Range("C2:F8").NumberFormat = "@"

$ws.Range("C2").Value = "9"
$ws.Range("D2").Value = "6"
$ws.Range("E2").Value = "0"
$ws.Range("F2").Value = "1"

$ws.Range("C3").Value = "48"
$ws.Range("D3").Value = "25"
$ws.Range("E3").Value = "4"
$ws.Range("F3").Value = "3"

$ws.Range("C4").Value = "36"
$ws.Range("D4").Value = "26"
$ws.Range("E4").Value = "3"
$ws.Range("F4").Value = "1"

$ws.Range("C5").Value = "54"
$ws.Range("D5").Value = "31"
$ws.Range("E5").Value = "4"
$ws.Range("F5").Value = "3"

$ws.Range("C6").Value = "74"
$ws.Range("D6").Value = "32"
$ws.Range("E6").Value = "1"
$ws.Range("F6").Value = "9"

$ws.Range("C7").Value = "1"
$ws.Range("D7").Value = "4"
$ws.Range("E7").Value = "0"
$ws.Range("F7").Value = "0"

$ws.Range("C8").Value = "0"
$ws.Range("D8").Value = "3"
$ws.Range("E8").Value = "0"
$ws.Range("F8").Value = "0"

Write-Host "done"
